# Weekly refresh of the Cilantro price series: a new week's observation is
# inserted at the top of the data block (row 4) and every existing
# observation shifts down one row, with the last row splitting into two
# (old row 29 -> new rows 29 and 30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4; this shifts rows 4:29 down to 5:30 and
# carries the date-format style from column D along with it.
$ws.Rows(4).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(4, 3).Value = "Ñuble"
$ws.Cells.Item(4, 4).Value = 44630
$ws.Cells.Item(4, 5).Value = 16
$ws.Cells.Item(4, 6).Value = 100112040
$ws.Cells.Item(4, 7).Value = "Cilantro"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 550
$ws.Cells.Item(4, 12).Value = 600
$ws.Cells.Item(4, 13).Value = 575
$ws.Cells.Item(4, 14).Value = "`$/atado 0,5 a 1 kilo"
$ws.Cells.Item(4, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(4, 16).Value = 575
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"
